$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = "PF3"
$ws.Range("C26").Value = "PF5"
$ws.Range("C27").Value = "PF6"
$ws.Range("C28").Value = "PF7"
$ws.Range("C29").Value = "PF8"
$ws.Range("C30").Value = "PF9"
$ws.Range("C31").Value = "PF10"
$ws.Range("D25").Value = "adc3"
